$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Slit2'
$ws.Cells.Item(2, 3).Value = 'Sdc1'
$ws.Cells.Item(2, 4).Value = 'ECs'
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.01246433333333333
$ws.Cells.Item(2, 8).Value = 0.037393
$ws.Cells.Item(2, 9).Value = 0.0065371131913745
$ws.Cells.Item(2, 10).Value = 0.006537113191374499
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.482117666666667
$ws.Cells.Item(2, 14).Value = 4.446353
$ws.Cells.Item(2, 15).Value = 0.1401829251394648
$ws.Cells.Item(2, 16).Value = 0.1401829251394648
$ws.Cells.Item(2, 17).Value = 0.01847360863655555
$ws.Cells.Item(2, 18).Value = 0.166262477729
$ws.Cells.Item(2, 19).Value = 0.0009163916491346595
$ws.Cells.Item(2, 20).Value = 0.0009163916491346592

$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Slit2'
$ws.Cells.Item(3, 3).Value = 'Sdc1'
$ws.Cells.Item(3, 4).Value = 'FAPs'
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.01246433333333333
$ws.Cells.Item(3, 8).Value = 0.037393
$ws.Cells.Item(3, 9).Value = 0.0065371131913745
$ws.Cells.Item(3, 10).Value = 0.006537113191374499
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.365790333333333
$ws.Cells.Item(3, 14).Value = 7.097371
$ws.Cells.Item(3, 15).Value = 0.2237632116883227
$ws.Cells.Item(3, 16).Value = 0.2237632116883226
$ws.Cells.Item(3, 17).Value = 0.02948799931144444
$ws.Cells.Item(3, 18).Value = 0.2653919938029999
$ws.Cells.Item(3, 19).Value = 0.001462765442872059
$ws.Cells.Item(3, 20).Value = 0.001462765442872058

$ws.Cells.Item(4, 1).Value = 'ECs'
$ws.Cells.Item(4, 2).Value = 'Slit2'
$ws.Cells.Item(4, 3).Value = 'Sdc1'
$ws.Cells.Item(4, 4).Value = 'sCs'
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.01246433333333333
$ws.Cells.Item(4, 8).Value = 0.037393
$ws.Cells.Item(4, 9).Value = 0.0065371131913745
$ws.Cells.Item(4, 10).Value = 0.006537113191374499
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 6.724832333333333
$ws.Cells.Item(4, 14).Value = 20.174497
$ws.Cells.Item(4, 15).Value = 0.6360538631722126
$ws.Cells.Item(4, 16).Value = 0.6360538631722126
$ws.Cells.Item(4, 17).Value = 0.08382055181344443
$ws.Cells.Item(4, 18).Value = 0.7543849663209998
$ws.Cells.Item(4, 19).Value = 0.004157956099367782
$ws.Cells.Item(4, 20).Value = 0.004157956099367781

$ws.Cells.Item(5, 1).Value = 'FAPs'
$ws.Cells.Item(5, 2).Value = 'Slit2'
$ws.Cells.Item(5, 3).Value = 'Sdc1'
$ws.Cells.Item(5, 4).Value = 'ECs'
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.9943730000000001
$ws.Cells.Item(5, 8).Value = 2.983119
$ws.Cells.Item(5, 9).Value = 0.5215143627507798
$ws.Cells.Item(5, 10).Value = 0.5215143627507798
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.482117666666667
$ws.Cells.Item(5, 14).Value = 4.446353
$ws.Cells.Item(5, 15).Value = 0.1401829251394648
$ws.Cells.Item(5, 16).Value = 0.1401829251394648
$ws.Cells.Item(5, 17).Value = 1.473777790556333
$ws.Cells.Item(5, 18).Value = 13.264000115007
$ws.Cells.Item(5, 19).Value = 0.07310740887264826
$ws.Cells.Item(5, 20).Value = 0.07310740887264824

$ws.Cells.Item(6, 1).Value = 'FAPs'
$ws.Cells.Item(6, 2).Value = 'Slit2'
$ws.Cells.Item(6, 3).Value = 'Sdc1'
$ws.Cells.Item(6, 4).Value = 'FAPs'
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.9943730000000001
$ws.Cells.Item(6, 8).Value = 2.983119
$ws.Cells.Item(6, 9).Value = 0.5215143627507798
$ws.Cells.Item(6, 10).Value = 0.5215143627507798
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.365790333333333
$ws.Cells.Item(6, 14).Value = 7.097371
$ws.Cells.Item(6, 15).Value = 0.2237632116883227
$ws.Cells.Item(6, 16).Value = 0.2237632116883226
$ws.Cells.Item(6, 17).Value = 2.352478031127667
$ws.Cells.Item(6, 18).Value = 21.172302280149
$ws.Cells.Item(6, 19).Value = 0.1166957287507034
$ws.Cells.Item(6, 20).Value = 0.1166957287507034

$ws.Cells.Item(7, 1).Value = 'FAPs'
$ws.Cells.Item(7, 2).Value = 'Slit2'
$ws.Cells.Item(7, 3).Value = 'Sdc1'
$ws.Cells.Item(7, 4).Value = 'sCs'
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.9943730000000001
$ws.Cells.Item(7, 8).Value = 2.983119
$ws.Cells.Item(7, 9).Value = 0.5215143627507798
$ws.Cells.Item(7, 10).Value = 0.5215143627507798
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 6.724832333333333
$ws.Cells.Item(7, 14).Value = 20.174497
$ws.Cells.Item(7, 15).Value = 0.6360538631722126
$ws.Cells.Item(7, 16).Value = 0.6360538631722126
$ws.Cells.Item(7, 17).Value = 6.686991701793667
$ws.Cells.Item(7, 18).Value = 60.182925316143
$ws.Cells.Item(7, 19).Value = 0.3317112251274281
$ws.Cells.Item(7, 20).Value = 0.3317112251274281

$ws.Cells.Item(8, 1).Value = 'sCs'
$ws.Cells.Item(8, 2).Value = 'Slit2'
$ws.Cells.Item(8, 3).Value = 'Sdc1'
$ws.Cells.Item(8, 4).Value = 'ECs'
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.8998656666666666
$ws.Cells.Item(8, 8).Value = 2.699597
$ws.Cells.Item(8, 9).Value = 0.4719485240578458
$ws.Cells.Item(8, 10).Value = 0.4719485240578457
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.482117666666667
$ws.Cells.Item(8, 14).Value = 4.446353
$ws.Cells.Item(8, 15).Value = 0.1401829251394648
$ws.Cells.Item(8, 16).Value = 0.1401829251394648
$ws.Cells.Item(8, 17).Value = 1.333706802193444
$ws.Cells.Item(8, 18).Value = 12.003361219741
$ws.Cells.Item(8, 19).Value = 0.0661591246176819
$ws.Cells.Item(8, 20).Value = 0.06615912461768188

$ws.Cells.Item(9, 1).Value = 'sCs'
$ws.Cells.Item(9, 2).Value = 'Slit2'
$ws.Cells.Item(9, 3).Value = 'Sdc1'
$ws.Cells.Item(9, 4).Value = 'FAPs'
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.8998656666666666
$ws.Cells.Item(9, 8).Value = 2.699597
$ws.Cells.Item(9, 9).Value = 0.4719485240578458
$ws.Cells.Item(9, 10).Value = 0.4719485240578457
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.365790333333333
$ws.Cells.Item(9, 14).Value = 7.097371
$ws.Cells.Item(9, 15).Value = 0.2237632116883227
$ws.Cells.Item(9, 16).Value = 0.2237632116883226
$ws.Cells.Item(9, 17).Value = 2.128893495498555
$ws.Cells.Item(9, 18).Value = 19.160041459487
$ws.Cells.Item(9, 19).Value = 0.1056047174947472
$ws.Cells.Item(9, 20).Value = 0.1056047174947471

$ws.Cells.Item(10, 1).Value = 'sCs'
$ws.Cells.Item(10, 2).Value = 'Slit2'
$ws.Cells.Item(10, 3).Value = 'Sdc1'
$ws.Cells.Item(10, 4).Value = 'sCs'
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.8998656666666666
$ws.Cells.Item(10, 8).Value = 2.699597
$ws.Cells.Item(10, 9).Value = 0.4719485240578458
$ws.Cells.Item(10, 10).Value = 0.4719485240578457
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 6.724832333333333
$ws.Cells.Item(10, 14).Value = 20.174497
$ws.Cells.Item(10, 15).Value = 0.6360538631722126
$ws.Cells.Item(10, 16).Value = 0.6360538631722126
$ws.Cells.Item(10, 17).Value = 6.051445730856554
$ws.Cells.Item(10, 18).Value = 54.46301157770899
$ws.Cells.Item(10, 19).Value = 0.3001846819454167
$ws.Cells.Item(10, 20).Value = 0.3001846819454166
